# Updated LHJ Population Data File (date updated: 2025.05.21)
# Adds a "last updated" note, a "Data Sources:" label, and the list of
# DOF/Census data sources below the existing Column/Description table.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New footnote-style rows under the table (row 8 intentionally left blank,
# matching the gap that existed between the table and these new lines).
$ws.Range("B9").Value  = "LHJ Population dataset last updated May 21, 2025"
$ws.Range("B10").Value = "Data Sources:"
$ws.Range("B11").Value = "DOF P3 Vintage 2025 (released 2025.04.25)"
$ws.Range("B12").Value = "DOF Annual Intercensal 2010-2020 (released 2025.01.29) and 2000-2010 (released 2013.03.19)"
$ws.Range("B13").Value = "DOF E4 2020-2025 (released May 2025), DOF E4 2010-2020 (released May 2025), and DOF E4 2000-2010 (released 2012.11.09)"
$ws.Range("B14").Value = "DOF E6 2020-2024 (released December 2024), DOF E6 2010-2019 (released December 2021), and DOF E6 2000-2010 (released December 2011)"
$ws.Range("B15").Value = "Census Decennial 2000, 2010, and 2020"

# Italicize the new notes (creates the new italic font / cell style).
$ws.Range("B9:B15").Font.Italic = $true

# Leave the cursor where the author left it after typing the new lines.
$ws.Range("B17").Select()
